$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = "fe0931d3-be07-4bed-9100-e63753bb21fd.md"
$ws.Range("B2").Value = "e2e\fe0931d3-be07-4bed-9100-e63753bb21fd.md"
$ws.Hyperlinks.Item(1).TextToDisplay = "e2e\fe0931d3-be07-4bed-9100-e63753bb21fd.md"
$ws.Range("G2").Value = "2016-08-21 19:02:33"

# --- zh-cn sheet ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A2").Value = "fe0931d3-be07-4bed-9100-e63753bb21fd.md"
$ws.Hyperlinks.Item(1).TextToDisplay = "fe0931d3-be07-4bed-9100-e63753bb21fd.md"
$ws.Range("G2").Value = "fe0931d3-be07-4bed-9100-e63753bb21fd.bc3c562acb8c1a3e1631dc28ff37990b5101490b.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-21 19:02:28"
$ws.Hyperlinks.Item(2).Delete()
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "0001-01-01 00:00:00"

# --- de-de sheet ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A2").Value = "fe0931d3-be07-4bed-9100-e63753bb21fd.md"
$ws.Hyperlinks.Item(1).TextToDisplay = "fe0931d3-be07-4bed-9100-e63753bb21fd.md"
$ws.Range("G2").Value = "fe0931d3-be07-4bed-9100-e63753bb21fd.bc3c562acb8c1a3e1631dc28ff37990b5101490b.de-de.xlf"
$ws.Range("H2").Value = "2016-08-21 19:02:33"
$ws.Hyperlinks.Item(2).Delete()
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "0001-01-01 00:00:00"

# --- Column width changes on zh-cn and de-de sheets (cols I & J) ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Columns.Item(9).ColumnWidth = 18.6506053379604
$ws.Columns.Item(10).ColumnWidth = 21.7054770333426

$ws = $wb.Worksheets.Item("de-de")
$ws.Columns.Item(9).ColumnWidth = 18.6506053379604
$ws.Columns.Item(10).ColumnWidth = 21.7054770333426

Write-Host "done"
